$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the header row: rename C1, add D1 and E1
$ws.Range("C1").Value = "Frecuencia del primer armonico"
$ws.Range("D1").Value = "Frecuencia del segundo armonico"
$ws.Range("E1").Value = "Frecuencia tercer armonico"

# Copy the header formatting (bold, border, centered) from C1 onto the new D1/E1 cells
$ws.Range("C1").Copy()
$ws.Range("D1:E1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

$cVals = @(302.0822542214855, 168.2561307901906, 260.8460358978987, 299.7026759167493, 144.0789473684217, 193.4083283994478, 150.2946005988606, 205.161499756533, 168.7681862269646, 181.3523448813639, 164.0458988525288, 189.4514362570849, 183.6957842303482, 324.6581385090431, 166.5045710951181, 171.2262632764723, 163.9610389610389, 166.9185167964661, 461.9599578503685, 379.0322580645161, 150.4509128235213, 339.2268504644899, 140.4676971858898, 188.1450360763547, 152.1725769774448, 163.3448029604624, 400.9995002498754, 332.0909908752092, 249.4767685223942, 257.8676330345224, 286.0475549573803, 246.03836530442, 496.5117809661706, 260.3076362974425, 230.2755011355785, 225.1522655779399, 215.3812683563583, 212.5879106204929, 278.0943697287439, 261.7621899059022, 277.8639890460972, 173.34556949118, 294.7187264709046, 209.6717754916372, 199.3869119019055, 202.4386170035077, 234.4582593250443, 242.9629629629635)
$dVals = @(150.4195586864189, 335.8310626702996, 395.0572742851982, 149.8513379583746, 290.1315789473683, 386.8166567988947, 300.2028397565919, 409.0245090082776, 335.5965082444227, 362.2591066057703, 328.0917977050576, 380.0557210106635, 368.3635038269954, 164.0935156594624, 336.1213771639759, 337.0453813968461, 321.4285714285716, 331.5122631640129, 230.9799789251847, 250.5040322580644, 301.1951022802259, 169.0140845070418, 422.6714229092349, 373.732760982738, 304.0993178046829, 488.9956501980132, 275.4622688655672, 166.5595681788973, 497.2791963164505, 513.4465000953651, 572.4540152534773, 491.2427022518768, 246.939581413716, 518.066806225539, 918.3371185938577, 452.1785690382167, 427.7167409985859, 424.8197275883558, 557.1005242762712, 784.4311377245513, 553.1720675490642, 345.4675231977158, 147.5744863934601, 420.1833578276996, 395.0419832067173, 404.6545292578362, 471.0479573712255, 484.8801742919395)
$eVals = @(450.4299181601573, 502.7247956403271, 130.9641922972851, 449.5540138751239, 434.8684210526312, 580.2249851983424, 454.7474162078624, 816.1012822593739, 504.3646944713873, 732.093126879804, 492.1376965575864, 569.1228744355844, 549.1434819584492, 490.5161005734453, 501.0698307722232, 516.2536208561314, 482.1428571428573, 501.6854585609672, 354.0569020021076, 126.0080645161288, 451.3527384705631, 1353.311357506743, 281.8866428854535, 562.2431272262306, 455.5343863315102, 1302.343699279361, 133.9330334832584, 498.1364863128128, 1246.267615459746, 3089.834064466909, 1433.826828174069, 984.9874895746452, 738.712649730156, 1038.682078820424, 689.2465685790457, 1129.509403654374, 1074.295659741108, 1273.034808154544, 834.2831091862317, 1833.190761334474, 831.4011866727524, 520.036708473539, 442.7234591803808, 339.5619007628247, 597.0944955351197, 608.2066700072378, 942.095914742451, 972.8976034858388)

# Rows 2..49 map to array indices 0..47
for ($i = 0; $i -lt $cVals.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 3).Value = $cVals[$i]
    $ws.Cells.Item($row, 4).Value = $dVals[$i]
    $ws.Cells.Item($row, 5).Value = $eVals[$i]
}

Write-Host "Done updating sheet"
